$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (single-decimal numeric-looking strings).
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "26.213.74"
$ws.Cells.Item(2, 5).Value = "  -2.16%  "
$ws.Cells.Item(3, 4).Value = "1.669.53"
$ws.Cells.Item(3, 5).Value = "  -1.61%  "
$ws.Cells.Item(4, 5).Value = "  +0.13%  "
$ws.Cells.Item(5, 4).Value = "217.73"
$ws.Cells.Item(5, 5).Value = "  -1.27%  "
$ws.Cells.Item(6, 4).Value = "0.5114"
$ws.Cells.Item(6, 5).Value = "  -0.17%  "
$ws.Cells.Item(7, 5).Value = "  +0.19%  "
$ws.Cells.Item(8, 4).Value = "0.2653"
$ws.Cells.Item(8, 5).Value = "  +2.94%  "
$ws.Cells.Item(9, 4).Value = "0.06371"
$ws.Cells.Item(9, 5).Value = "  +2.91%  "
$ws.Cells.Item(10, 5).Value = "  -2.50%  "
$ws.Cells.Item(11, 4).Value = "0.07386"
$ws.Cells.Item(11, 5).Value = "  +0.41%  "
$ws.Cells.Item(12, 5).Value = "  +1.33%  "
$ws.Cells.Item(13, 4).Value = "1.672.81"
$ws.Cells.Item(13, 5).Value = "  -1.40%  "
$ws.Cells.Item(14, 4).Value = "0.5814"
$ws.Cells.Item(14, 5).Value = "  -0.04%  "
$ws.Cells.Item(15, 4).Value = "0.000008629"
$ws.Cells.Item(15, 5).Value = "  +5.31%  "
$ws.Cells.Item(16, 4).Value = "64.37"
$ws.Cells.Item(16, 5).Value = "  -1.47%  "
$ws.Cells.Item(17, 4).Value = "26.266.06"
$ws.Cells.Item(17, 5).Value = "  -2.04%  "
$ws.Cells.Item(18, 4).Value = "4.931"
$ws.Cells.Item(18, 5).Value = "  -1.88%  "
$ws.Cells.Item(19, 4).Value = "1.007"
$ws.Cells.Item(19, 5).Value = "  +0.14%  "
$ws.Cells.Item(20, 5).Value = "  +1.83%  "
$ws.Cells.Item(21, 4).Value = "188.83"
$ws.Cells.Item(21, 5).Value = "  +0.82%  "
$ws.Cells.Item(22, 4).Value = "6.201"
$ws.Cells.Item(22, 5).Value = "  -1.14%  "
$ws.Cells.Item(23, 5).Value = "  +0.18%  "
$ws.Cells.Item(24, 4).Value = "144.22"
$ws.Cells.Item(25, 4).Value = "7.634"
$ws.Cells.Item(25, 5).Value = "  +1.83%  "
$ws.Cells.Item(26, 4).Value = "0.1175"
$ws.Cells.Item(26, 5).Value = "  +2.45%  "
$ws.Cells.Item(27, 4).Value = "15.62"
$ws.Cells.Item(27, 5).Value = "  +2.66%  "
$ws.Cells.Item(28, 4).Value = "0.05962"
$ws.Cells.Item(28, 5).Value = "  +1.32%  "
$ws.Cells.Item(29, 5).Value = "  -3.74%  "
$ws.Cells.Item(30, 5).Value = "  -2.00%  "
$ws.Cells.Item(31, 4).Value = "3.519"
$ws.Cells.Item(31, 5).Value = "  +1.59%  "
$ws.Cells.Item(32, 4).Value = "3.520"
$ws.Cells.Item(32, 5).Value = "  +2.23%  "
$ws.Cells.Item(33, 4).Value = "1.642"
$ws.Cells.Item(33, 5).Value = "  -0.36%  "
$ws.Cells.Item(34, 4).Value = "1.012"
$ws.Cells.Item(34, 5).Value = "  +2.02%  "
$ws.Cells.Item(35, 4).Value = "0.6025"
$ws.Cells.Item(35, 5).Value = "  +0.85%  "
$ws.Cells.Item(36, 4).Value = "2.377"
$ws.Cells.Item(36, 5).Value = "  -1.61%  "
$ws.Cells.Item(37, 4).Value = "2.657"
$ws.Cells.Item(37, 5).Value = "  -0.48%  "
$ws.Cells.Item(38, 4).Value = "6.092"
$ws.Cells.Item(38, 5).Value = "  +4.25%  "
$ws.Cells.Item(39, 4).Value = "0.01615"
$ws.Cells.Item(39, 5).Value = "  +1.08%  "
$ws.Cells.Item(40, 4).Value = "1.077.66"
$ws.Cells.Item(40, 5).Value = "  -1.33%  "
$ws.Cells.Item(41, 5).Value = "  +1.28%  "
$ws.Cells.Item(42, 5).Value = "  +0.55%  "
$ws.Cells.Item(43, 4).Value = "100.26"
$ws.Cells.Item(43, 5).Value = "  +3.05%  "
$ws.Cells.Item(44, 4).Value = "1.820.92"
$ws.Cells.Item(44, 5).Value = "  -1.11%  "
$ws.Cells.Item(45, 4).Value = "0.00000000114"
$ws.Cells.Item(45, 5).Value = "  +9.82%  "
$ws.Cells.Item(46, 4).Value = "56.19"
$ws.Cells.Item(46, 5).Value = "  +0.19%  "
$ws.Cells.Item(47, 5).Value = "  +1.05%  "
$ws.Cells.Item(48, 4).Value = "8.040"
$ws.Cells.Item(48, 5).Value = "  +0.21%  "
$ws.Cells.Item(49, 4).Value = "0.05210"
$ws.Cells.Item(49, 5).Value = "  -0.41%  "
$ws.Cells.Item(50, 5).Value = "  -0.54%  "
$ws.Cells.Item(51, 4).Value = "5.874"
$ws.Cells.Item(51, 5).Value = "  +2.17%  "
